$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2167615066290701
$ws.Range("C2").Value = 5018056568004490

$ws.Range("B3").Value = 2137973970787784
$ws.Range("C3").Value = 4949437397367131

$ws.Range("B4").Value = 2132527638753326
$ws.Range("C4").Value = 4936830349835802

$ws.Range("B5").Value = 1543741264872729
$ws.Range("C5").Value = 3573784987446494

$ws.Range("B6").Value = 944229903108372.4
$ws.Range("C6").Value = 2185908177301574

$ws.Range("B7").Value = 804472081958416
$ws.Range("C7").Value = 1862362766069293

$ws.Range("B8").Value = 666816676091205.1
$ws.Range("C8").Value = 1543692196125669

$ws.Range("B9").Value = 109414336465166.3
$ws.Range("C9").Value = 253296144801917.1

$ws.Range("B10").Value = 112189862.6265715
$ws.Range("C10").Value = 259720789.6349103
